$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "Prajatantra diwas"
$ws.Range("C10").Value = "Class and Object"

$ws.Rows.Item(20).RowHeight = 20.25
